# Auto-generated edit script applying the cryptos.xlsx diff (price + volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.492.65"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.921.65"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.01"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2876"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06754"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "106.54"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.29"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07763"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").Value = "1.913.27"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.300"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6608"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "292.19"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.74%  "
$ws.Range("D17").Value = "30.499.27"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007589"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.93"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "2.148.08"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.291"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.35%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9997"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.210"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.368"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.29"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.30"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.098"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1070"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.367"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.173"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05037"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7436"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02114"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.683"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.074"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.25"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8746"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.930"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "67.52"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "49.68"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +17.97%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.193"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.316"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.15"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2449"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.80%  "
